$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows starting at row 864; this shifts every
# existing row from 864 onward down by two rows (864->866, ..., 948->950)
$ws.Rows.Item(864).Resize(2).Insert()

# Fill in the first new row (864) with its data
$ws.Range("A864").Value = 4
$ws.Range("B864").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C864").Value = "Los Lagos"
$ws.Range("D864").Value = "2023-08-28"
$ws.Range("E864").Value = 10
$ws.Range("F864").Value = "Fruta"
$ws.Range("G864").Value = 100108
$ws.Range("H864").Value = "Tropicales y subtropicales"
$ws.Range("I864").Value = 100108006
$ws.Range("J864").Value = "Plátano"
$ws.Range("K864").Value = "Sin especificar"
$ws.Range("L864").Value = "Pintón"
$ws.Range("M864").Value = 200
$ws.Range("N864").Value = 18000
$ws.Range("O864").Value = 18000
$ws.Range("P864").Value = 18000
$ws.Range("Q864").Value = "$/caja 20 kilos"
$ws.Range("R864").Value = "Ecuador"
$ws.Range("S864").Value = 900
$ws.Range("T864").Value = 20

# Fill in the second new row (865) with its data
$ws.Range("A865").Value = 4
$ws.Range("B865").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C865").Value = "Los Lagos"
$ws.Range("D865").Value = "2023-08-28"
$ws.Range("E865").Value = 10
$ws.Range("F865").Value = "Fruta"
$ws.Range("G865").Value = 100108
$ws.Range("H865").Value = "Tropicales y subtropicales"
$ws.Range("I865").Value = 100108006
$ws.Range("J865").Value = "Plátano"
$ws.Range("K865").Value = "Sin especificar"
$ws.Range("L865").Value = "Primera Pintón"
$ws.Range("M865").Value = 300
$ws.Range("N865").Value = 20000
$ws.Range("O865").Value = 20000
$ws.Range("P865").Value = 20000
$ws.Range("Q865").Value = "$/caja 20 kilos"
$ws.Range("R865").Value = "Ecuador"
$ws.Range("S865").Value = 1000
$ws.Range("T865").Value = 20
